# Threat Alert Report refresh (2026-01-30 06:29)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New tracking row (row 5): clone row 4's layout/format first so the new
#     row's cells pick up the same borders/fill/alignment as the rest of the
#     table before we touch any values. ---
$ws.Range("A4:K4").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 2: new market check ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "13-FEB-26"
$ws.Range("D2").Value = 687
$ws.Range("F2").Value = -208

# --- Row 4: new market check ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "13-MAR-26"
$ws.Range("D4").Value = 1237
$ws.Range("E4").Value = 1501
$ws.Range("F4").Value = -264

# --- Row 5: carries forward what used to be row 4's reading ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "27-MAR-26"
$ws.Range("B5").Value = "SM-446"
$ws.Range("C5").Value = "Air Arabia Egypt E5-512"
$ws.Range("D5").Value = 513
$ws.Range("E5").Value = 786
$ws.Range("F5").Value = -273
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"
